# This edit re-orders the data rows (2-10) of the "Artfynd" sheet.
# Each full row of data (identified by its unique Id in column A) is moved
# to a different row position; the content of every row stays identical,
# only its row position changes (a permutation of rows 2..10).
#
# Mapping is expressed as: sourceRow (current position) -> destRow (new position)
#   2 -> 7, 3 -> 4, 4 -> 10, 5 -> 6, 6 -> 9, 7 -> 5, 8 -> 2, 9 -> 3, 10 -> 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 10
$firstCol = 1
$lastCol = 51   # column AY

$mapping = @{
    2  = 7
    3  = 4
    4  = 10
    5  = 6
    6  = 9
    7  = 5
    8  = 2
    9  = 3
    10 = 8
}

# Some columns store values that *look* like numbers/dates/times but are
# actually plain TEXT in the source file:
#   I  (9)  "Antal"     - small numeric-looking text values ("30","4",...)
#   Y  (25) "Startdatum" - date-looking text "2023-07-17"
#   Z  (26) "Starttid"   - time-looking text "00:00"
#   AA (27) "Slutdatum"  - date-looking text "2023-07-17"
#   AB (28) "Sluttid"    - time-looking text "00:00"
# Track those so we can re-apply text formatting (quote-prefix trick) when
# moving the values to their new row, preventing Excel from "helpfully"
# re-interpreting them as numbers/dates.
$textColumns = @(9, 25, 26, 27, 28)

# --- Step 1: buffer every source cell's value (and type) into memory ---
$buffer = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        $key = "$r-$c"
        $buffer[$key] = $val
    }
}

# --- Step 2: clear out all the cells in the affected range ---
$clearRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$clearRange.ClearContents()

# --- Step 3: write the buffered values to their new (destination) row ---
for ($srcRow = $firstRow; $srcRow -le $lastRow; $srcRow++) {
    $destRow = $mapping[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $key = "$srcRow-$c"
        $val = $buffer[$key]

        if ($val -eq $null) {
            continue
        }
        if ($val -is [string] -and $val -eq "") {
            continue
        }

        $destCell = $ws.Cells.Item($destRow, $c)

        if ($textColumns -contains $c) {
            # Force the value to be stored as text, even though it looks
            # like a number (e.g. "30"), matching the original file.
            $destCell.Value = "'" + $val
        } else {
            $destCell.Value = $val
        }
    }
}

Write-Output "Row reorder complete"
